$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.870.94'
$ws.Range("E2").Value = '  +0.07%  '

# Row 3
$ws.Range("D3").Value = '2.572.54'
$ws.Range("E3").Value = '  +1.60%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").Value = "'313.29"
$ws.Range("E5").Value = '  -0.71%  '

# Row 6
$ws.Range("E6").Value = '  +3.63%  '

# Row 7
$ws.Range("D7").Value = "'0.574"
$ws.Range("E7").Value = '  -0.41%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = '  +0.14%  '

# Row 10
$ws.Range("D10").Value = "'35.87"
$ws.Range("E10").Value = '  -0.84%  '

# Row 11
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = '  +0.45%  '

# Row 12
$ws.Range("D12").Value = "'7.45"
$ws.Range("E12").Value = '  -1.43%  '

# Row 13
$ws.Range("D13").Value = '2.962.40'
$ws.Range("E13").Value = '  +1.46%  '

# Row 14
$ws.Range("E14").Value = '  -1.23%  '

# Row 15
$ws.Range("D15").Value = "'15.95"
$ws.Range("E15").Value = '  +4.78%  '

# Row 16
$ws.Range("D16").Value = '2.611.49'
$ws.Range("E16").Value = '  +0.19%  '

# Row 17
$ws.Range("D17").Value = "'0.847"
$ws.Range("E17").Value = '  -0.59%  '

# Row 18
$ws.Range("D18").Value = '42.895.84'
$ws.Range("E18").Value = '  +0.02%  '

# Row 19
$ws.Range("E19").Value = '  -0.93%  '

# Row 20
$ws.Range("D20").Value = "'12.59"
$ws.Range("E20").Value = '  -3.91%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0964'
$ws.Range("E21").Value = '  -0.10%  '

# Row 22
$ws.Range("E22").Value = '  -0.59%  '

# Row 23
$ws.Range("D23").Value = "'250.31"
$ws.Range("E23").Value = '  -1.20%  '

# Row 24
$ws.Range("E24").Value = '  +0.41%  '

# Row 25
$ws.Range("E25").Value = '  +0.04%  '

# Row 26
$ws.Range("D26").Value = "'27.18"
$ws.Range("E26").Value = '  +2.03%  '

# Row 27
$ws.Range("E27").Value = '  -0.02%  '

# Row 28
$ws.Range("E28").Value = '  -0.71%  '

# Row 29
$ws.Range("D29").Value = "'39.88"
$ws.Range("E29").Value = '  -1.59%  '

# Row 30
$ws.Range("D30").Value = "'10.28"
$ws.Range("E30").Value = '  -1.43%  '

# Row 31
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = "'158.16"
$ws.Range("E31").Value = '  +0.10%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'5.81"
$ws.Range("E32").Value = '  -2.08%  '

# Row 33
$ws.Range("E33").Value = '  +1.10%  '

# Row 34
$ws.Range("E34").Value = '  +2.56%  '

# Row 35
$ws.Range("E35").Value = '  -2.17%  '

# Row 36
$ws.Range("E36").Value = '  +0.19%  '

# Row 37
$ws.Range("D37").Value = "'18.68"
$ws.Range("E37").Value = '  -2.32%  '

# Row 38
$ws.Range("E38").Value = '  +11.16%  '

# Row 39
$ws.Range("E39").Value = '  +0.21%  '

# Row 40
$ws.Range("E40").Value = '  -0.16%  '

# Row 41
$ws.Range("D41").Value = "'23.37"
$ws.Range("E41").Value = '  +0.21%  '

# Row 42
$ws.Range("D42").Value = "'4.13"
$ws.Range("E42").Value = '  +7.36%  '

# Row 43
$ws.Range("E43").Value = '  -0.59%  '

# Row 44
$ws.Range("E44").Value = '  -0.13%  '

# Row 45
$ws.Range("D45").Value = "'3.25"
$ws.Range("E45").Value = '  -1.89%  '

# Row 46
$ws.Range("D46").Value = '1.999.38'
$ws.Range("E46").Value = '  -1.49%  '

# Row 47
$ws.Range("E47").Value = '  -2.11%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = "'0.198"
$ws.Range("E48").Value = '  +3.05%  '

# Row 49
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.814.25'
$ws.Range("E49").Value = '  +1.46%  '

# Row 50
$ws.Range("D50").Value = "'82.06"
$ws.Range("E50").Value = '  -3.48%  '

# Row 51
$ws.Range("D51").Value = "'74.73"
$ws.Range("E51").Value = '  -0.37%  '
